# Natmi following Dr Hou advice
# Update the LR-pair statistics sheet: the number of ligand-/receptor-expressing
# cells changed from 1 to 3 per group, which changes average/total expression,
# specificity scores and edge weights for every data row (rows 2-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New count of expressing cells (both ligand side "E" and receptor side "K")
$newCount = 3

# Ligand average expression value (column G) keyed by "Sending cluster" (column A)
$gBySending = @{
    "ECs"    = 4.029909
    "FAPs"   = 3.678840333333333
    "M1"     = 9.448480333333332
    "M2"     = 11.30319166666666
    "Neutro" = 3.740567333333333
    "sCs"    = 2.173862333333334
}

# Receptor average expression value (column M) keyed by "Target cluster" (column D)
$mByTarget = @{
    "ECs"    = 3.175664333333333
    "FAPs"   = 19.658112
    "M1"     = 9.985787333333334
    "M2"     = 10.739333
    "Neutro" = 8.480710666666667
    "sCs"    = 23.366118
}

$gSum = 0
foreach ($v in $gBySending.Values) { $gSum += $v }

$mSum = 0
foreach ($v in $mByTarget.Values) { $mSum += $v }

$firstRow = 2
$lastRow = 37

# First pass: compute per-row Q/R so the global sums used for S/T (specificity
# of the edge weight) can be derived.
$qVals = @{}
$rVals = @{}
$qSum = 0
$rSum = 0

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $sending = [string]$ws.Cells.Item($row, 1).Value2   # column A
    $target  = [string]$ws.Cells.Item($row, 4).Value2   # column D

    $g = $gBySending[$sending]
    $m = $mByTarget[$target]

    $h = $g * $newCount
    $n = $m * $newCount

    $q = $g * $m
    $r = $h * $n

    $qVals[$row] = $q
    $rVals[$row] = $r

    $qSum += $q
    $rSum += $r
}

# Second pass: write all the updated values.
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $sending = [string]$ws.Cells.Item($row, 1).Value2   # column A
    $target  = [string]$ws.Cells.Item($row, 4).Value2   # column D

    $g = $gBySending[$sending]
    $m = $mByTarget[$target]

    $h = $g * $newCount
    $n = $m * $newCount

    $i = $g / $gSum
    $o = $m / $mSum

    $q = $qVals[$row]
    $r = $rVals[$row]

    $s = $q / $qSum
    $t = $r / $rSum

    $ws.Cells.Item($row, 5).Value2  = $newCount   # E - Ligand-expressing cells
    $ws.Cells.Item($row, 7).Value2  = $g           # G - Ligand average expression value
    $ws.Cells.Item($row, 8).Value2  = $h           # H - Ligand total expression value
    $ws.Cells.Item($row, 9).Value2  = $i           # I - Ligand specificity (avg)
    $ws.Cells.Item($row, 10).Value2 = $i           # J - Ligand specificity (total)
    $ws.Cells.Item($row, 11).Value2 = $newCount    # K - Receptor-expressing cells
    $ws.Cells.Item($row, 13).Value2 = $m           # M - Receptor average expression value
    $ws.Cells.Item($row, 14).Value2 = $n           # N - Receptor total expression value
    $ws.Cells.Item($row, 15).Value2 = $o           # O - Receptor specificity (avg)
    $ws.Cells.Item($row, 16).Value2 = $o           # P - Receptor specificity (total)
    $ws.Cells.Item($row, 17).Value2 = $q           # Q - Edge average expression weight
    $ws.Cells.Item($row, 18).Value2 = $r           # R - Edge total expression weight
    $ws.Cells.Item($row, 19).Value2 = $s           # S - Edge specificity (avg)
    $ws.Cells.Item($row, 20).Value2 = $t           # T - Edge specificity (total)
}
